$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cell values for rows 2-7 (columns A:T) reflecting the new TPM data.
# Each inner array is ordered: A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$rows = @(2,3,4,5,6,7)
$values = @(
    @("ECs","Fgf16","Fgfr3","ECs",2,1,0.2355055,0.471011,0.2095102068231121,0.2095102068231121,2,1,7.6704545,15.340909,0.8278663930876066,0.7913005936208135,1.80643422224975,7.225736888999,0.1734464592376883,0.165785551028748),
    @("ECs","Fgf16","Fgfr3","FAPs",2,1,0.2355055,0.471011,0.2095102068231121,0.2095102068231121,3,1,0.856297,2.568891,0.0924194920655273,0.1325061620042962,0.2016626531335,1.209975918801,0.01936282689713559,0.0277613934068569),
    @("ECs","Fgf16","Fgfr3","MuSCs",2,1,0.2355055,0.471011,0.2095102068231121,0.2095102068231121,2,1,0.7385775,1.477155,0.0797141148468662,0.07619324437489022,0.17393906342625,0.695756253705,0.01670092068828825,0.01596326238750717),
    @("MuSCs","Fgf16","Fgfr3","ECs",2,1,0.888571,1.777142,0.7904897931768879,0.7904897931768879,2,1,7.6704545,15.340909,0.8278663930876066,0.7913005936208135,6.8157434255195,27.262973702078,0.6544199338499184,0.6255150425920655),
    @("MuSCs","Fgf16","Fgfr3","FAPs",2,1,0.888571,1.777142,0.7904897931768879,0.7904897931768879,3,1,0.856297,2.568891,0.0924194920655273,0.1325061620042962,0.760880681587,4.565284089522,0.07305666516839171,0.1047447685974393),
    @("MuSCs","Fgf16","Fgfr3","MuSCs",2,1,0.888571,1.777142,0.7904897931768879,0.7904897931768879,2,1,0.7385775,1.477155,0.0797141148468662,0.07619324437489022,0.6562785477525001,2.62511419101,0.06301319415857796,0.06022998198738305)
)

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $r = $rows[$ri]
    $vals = $values[$ri]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
